$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 101
$ws.Range("I5").Value = 101
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 101
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 14
$ws.Range("N5").ClearContents()
$ws.Range("H40").Value = 2748.7
$ws.Range("I40").Value = 2751
$ws.Range("J40").Value = 2746.4
$ws.Range("K40").Value = 2751
$ws.Range("L40").Value = 2746.4
$ws.Range("M40").Value = -2576
$ws.Range("N40").Value = -3096.4
$ws.Range("H74").Value = 4622.5
$ws.Range("I74").Value = 4660
$ws.Range("J74").Value = 4600
$ws.Range("K74").Value = 4660
$ws.Range("L74").Value = 4600
$ws.Range("M74").Value = -3724
$ws.Range("N74").Value = -6472
$ws.Range("H76").Value = 3334.1228
$ws.Range("I76").Value = 2982.7778
$ws.Range("J76").Value = 4651.6665
$ws.Range("K76").Value = 2982.7778
$ws.Range("L76").Value = 4651.6665
$ws.Range("M76").Value = -2667.7778
$ws.Range("N76").Value = -5281.6665
$ws.Range("H77").Value = 4622.5
$ws.Range("I77").Value = 4660
$ws.Range("J77").Value = 4600
$ws.Range("K77").Value = 23300
$ws.Range("L77").Value = 23000
$ws.Range("M77").Value = -18620
$ws.Range("N77").Value = -32360
$ws.Range("H79").Value = 3334.1228
$ws.Range("I79").Value = 2982.7778
$ws.Range("J79").Value = 4651.6665
$ws.Range("K79").Value = 2982.7778
$ws.Range("L79").Value = 4651.6665
$ws.Range("M79").Value = -1890.7778
$ws.Range("N79").Value = -6835.6665

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2578.3333
$ws.Range("I63").Value = 1867.5
$ws.Range("J63").Value = 4000
$ws.Range("K63").Value = 1867.5
$ws.Range("L63").Value = 4000
$ws.Range("M63").Value = -1181.5
$ws.Range("N63").Value = -5372
$ws.Range("H66").Value = 2578.3333
$ws.Range("I66").Value = 1867.5
$ws.Range("J66").Value = 4000
$ws.Range("K66").Value = 9337.5
$ws.Range("L66").Value = 20000
$ws.Range("M66").Value = -5905.5
$ws.Range("N66").Value = -26864
$ws.Range("H102").Value = 2493.3333
$ws.Range("I102").Value = 2493.3333
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 2493.3333
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -871.3332999999998
$ws.Range("N102").ClearContents()
$ws.Range("H122").Value = 2449
$ws.Range("I122").Value = 1781.6666
$ws.Range("K122").Value = 5344.9998
$ws.Range("M122").Value = -2894.9998

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2679
$ws.Range("I105").Value = 2342.8572
$ws.Range("J105").Value = 3463.3333
$ws.Range("K105").Value = 2342.8572
$ws.Range("L105").Value = 3463.3333
$ws.Range("M105").Value = -595.8571999999999
$ws.Range("N105").Value = -6957.3333

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3923.0833
$ws.Range("I62").Value = 2121.1667
$ws.Range("J62").Value = 5725
$ws.Range("K62").Value = 2121.1667
$ws.Range("L62").Value = 5725
$ws.Range("M62").Value = -1497.1667
$ws.Range("N62").Value = -6973
$ws.Range("H65").Value = 3923.0833
$ws.Range("I65").Value = 2121.1667
$ws.Range("J65").Value = 5725
$ws.Range("K65").Value = 10605.8335
$ws.Range("L65").Value = 28625
$ws.Range("M65").Value = -7485.833500000001
$ws.Range("N65").Value = -34865
$ws.Range("H122").Value = 460846.12
$ws.Range("I122").Value = 68031.13
$ws.Range("J122").Value = 1115537.8
$ws.Range("K122").Value = 204093.39
$ws.Range("L122").Value = 3346613.4
$ws.Range("M122").Value = -201643.39
$ws.Range("N122").Value = -3351513.4

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 585832.6
$ws.Range("I5").Value = 456
$ws.Range("J5").Value = 1254834.4
$ws.Range("K5").Value = 1368
$ws.Range("L5").Value = 3764503.2
$ws.Range("M5").Value = -1256
$ws.Range("N5").Value = -3764727.2
$ws.Range("H122").Value = 1101.5
$ws.Range("J122").Value = 3111.8
$ws.Range("L122").Value = 28006.2
$ws.Range("N122").Value = -32906.2
$ws.Range("H135").Value = 585832.6
$ws.Range("I135").Value = 456
$ws.Range("J135").Value = 1254834.4
$ws.Range("K135").Value = 4104
$ws.Range("L135").Value = 11293509.6
$ws.Range("M135").Value = -1569
$ws.Range("N135").Value = -11298579.6

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5408.306
$ws.Range("I70").Value = 4856.7
$ws.Range("J70").Value = 5788.724
$ws.Range("K70").Value = 4856.7
$ws.Range("L70").Value = 5788.724
$ws.Range("M70").Value = -4586.7
$ws.Range("N70").Value = -6328.724
$ws.Range("H73").Value = 5408.306
$ws.Range("I73").Value = 4856.7
$ws.Range("J73").Value = 5788.724
$ws.Range("K73").Value = 4856.7
$ws.Range("L73").Value = 5788.724
$ws.Range("M73").Value = -3920.7
$ws.Range("N73").Value = -7660.724
$ws.Range("H80").Value = 2564.4443
$ws.Range("I80").Value = 2593.3333
$ws.Range("J80").Value = 2550
$ws.Range("K80").Value = 2593.3333
$ws.Range("L80").Value = 2550
$ws.Range("M80").Value = -1595.3333
$ws.Range("N80").Value = -4546
$ws.Range("H83").Value = 2564.4443
$ws.Range("I83").Value = 2593.3333
$ws.Range("J83").Value = 2550
$ws.Range("K83").Value = 12966.6665
$ws.Range("L83").Value = 12750
$ws.Range("M83").Value = -7974.666499999999
$ws.Range("N83").Value = -22734
$ws.Range("H97").Value = 1567.8572
$ws.Range("I97").Value = 1645.8334
$ws.Range("J97").Value = 1100
$ws.Range("K97").Value = 1645.8334
$ws.Range("L97").Value = 1100
$ws.Range("M97").Value = -1149.8334
$ws.Range("N97").Value = -2092
$ws.Range("H122").Value = 121568.14
$ws.Range("I122").Value = 209714.67
$ws.Range("J122").Value = 4039.4443
$ws.Range("K122").Value = 629144.01
$ws.Range("L122").Value = 12118.3329
$ws.Range("M122").Value = -626694.01
$ws.Range("N122").Value = -17018.3329

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2452.0908
$ws.Range("I100").Value = 1271.125
$ws.Range("J100").Value = 5601.3335
$ws.Range("K100").Value = 1271.125
$ws.Range("L100").Value = 5601.3335
$ws.Range("M100").Value = -730.125
$ws.Range("N100").Value = -6683.3335
$ws.Range("H122").Value = 10103470
$ws.Range("I122").Value = 22223762
$ws.Range("J122").Value = 3226.6667
$ws.Range("K122").Value = 66671286
$ws.Range("L122").Value = 9680.000100000001
$ws.Range("M122").Value = -66668836
$ws.Range("N122").Value = -14580.0001

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 900.2857
$ws.Range("I100").Value = 900.4
$ws.Range("K100").Value = 1800.8
$ws.Range("M100").Value = -1259.8
$ws.Range("H122").Value = 31870.121
$ws.Range("I122").Value = 41381.56
$ws.Range("K122").Value = 124144.68
$ws.Range("M122").Value = -121694.68
